$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new trade record as row 5
$ws.Cells.Item(5, 1).Value = 9959.23
$ws.Cells.Item(5, 2).Value = 10035.5
$ws.Cells.Item(5, 3).Value = 109.08
$ws.Cells.Item(5, 4).Value = 108.25
$ws.Cells.Item(5, 5).Value = $false
$ws.Cells.Item(5, 6).Value = -0.76
$ws.Cells.Item(5, 7).Value = 42612.674525462964
$ws.Cells.Item(5, 8).Value = $false

# Match the date formatting used by the rows above (column G) by copying
# the existing cell's style instead of fabricating a new number format.
$ws.Cells.Item(4, 7).Copy()
$ws.Cells.Item(5, 7).PasteSpecial(-4122)
